$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# E2 already carries the plain data-row style (no number format / border /
# quote-prefix) that every new data cell in columns A, B and E should end
# up with. Copy-pasting its *formats only* (xlPasteFormats = -4122) lets us
# stamp that style onto a cell without touching its value - which is how
# the diff's trailing "<c r=.. s="1"/>" placeholder cells (no t=, no v=)
# get created, and also how we undo the auto quote-prefix style Excel
# applies when a cell value is forced to text with a leading apostrophe.
$blankStyleSource = $ws.Range("E2")

function Set-EmptyTextCell($addr) {
    $cell = $ws.Range($addr)
    # A literal leading apostrophe forces a genuine empty-string text value
    # (t="s" pointing at an empty shared string) instead of clearing the
    # cell outright, which is what a plain "" assignment would do.
    $cell.Value = "'"
    $blankStyleSource.Copy()
    $cell.PasteSpecial(-4122)
}

function Set-StyledBlankCell($addr) {
    $blankStyleSource.Copy()
    $ws.Range($addr).PasteSpecial(-4122)
}

# Row 2
Set-EmptyTextCell "A2"
Set-EmptyTextCell "B2"
$ws.Range("C2").Value = 257
$ws.Range("D2").Value = "20/05/2025"
Set-StyledBlankCell "E2"

# Row 3
Set-EmptyTextCell "A3"
Set-EmptyTextCell "B3"
$ws.Range("C3").Value = 227
$ws.Range("D3").Value = "20/05/2025"
Set-StyledBlankCell "E3"

# Row 4
Set-EmptyTextCell "A4"
Set-EmptyTextCell "B4"
$ws.Range("C4").Value = 115
$ws.Range("D4").Value = "20/05/2025"
Set-StyledBlankCell "E4"

# Row 5
Set-EmptyTextCell "A5"
Set-EmptyTextCell "B5"
$ws.Range("C5").Value = 102
$ws.Range("D5").Value = "20/05/2025"
Set-StyledBlankCell "E5"

# Row 6
$ws.Range("A6").Value = "S"
$ws.Range("B6").Value = "s@d.com"
$ws.Range("C6").Value = 96
$ws.Range("D6").Value = "19/05/2025"
Set-StyledBlankCell "E6"

# Row 7
$ws.Range("A7").Value = "ssssssss"
$ws.Range("B7").Value = "s"
$ws.Range("C7").Value = 49
$ws.Range("D7").Value = "20/05/2025"
Set-StyledBlankCell "E7"
